$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the "Test Set (81 Compounds)" label cell and repeat it on every row ---
$testSetLabel = $ws.Range("A3").Text
$ws.Range("A3:A9").UnMerge()
$ws.Range("A4:A9").Value = $testSetLabel

# --- Remove the horizontal-centering from column A (rows 3-9); keep vertical-center + wrap ---
$ws.Range("A3:A9").HorizontalAlignment = 1

# --- Tighten the numeric display format for the metric columns (10 decimals -> 5 decimals) ---
$ws.Range("C3:D9").NumberFormat = "0.00000"

# --- Update row heights: header row 2 shrinks, data rows 4-9 grow to match row 3 ---
$ws.Rows(2).RowHeight = 31.2
$ws.Rows(3).RowHeight = 15.6
$ws.Rows(4).RowHeight = 31.2
$ws.Rows(5).RowHeight = 31.2
$ws.Rows(6).RowHeight = 31.2
$ws.Rows(7).RowHeight = 31.2
$ws.Rows(8).RowHeight = 31.2
$ws.Rows(9).RowHeight = 31.2

# --- Move the active selection like in the saved file ---
$ws.Range("F9").Select()
